# Fixed naive component forecaster bug - Presentation state 11.02.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.1103015140893087, 1.867466406311011, 16.07388988832066, 4.009225597084886, 4.097780106789175, 23)
    3  = @(0.1241568715908939, 1.713852921858244, 10.94170538271163, 3.307824871832188, 3.383280951883025, 22)
    4  = @(-0.5207749427252132, 1.008444475967709, 4.236175762135406, 2.058197211672245, 2.040396452217854, 21)
    5  = @(0.08713849139876084, 0.7045681055244151, 1.740827324384067, 1.319404155057906, 1.350724642826101, 20)
    6  = @(0.05311708210651798, 0.7222531657823513, 1.611489849032536, 1.269444701053392, 1.30308821368913, 19)
    7  = @(0.01086339460369382, 0.5852842659497787, 0.6191602100143899, 0.7868673395270577, 0.8096026454595719, 18)
    8  = @(0.09509580837041771, 0.5696764806980906, 0.5776344400534255, 0.760022657592144, 0.7772568294901242, 17)
    9  = @(0.2285112237303573, 0.484995177258808, 0.3855122973504174, 0.6208963660309322, 0.5962504332730643, 16)
    10 = @(0.2058289907268505, 0.4532524315174537, 0.3699707508218096, 0.6082522098782788, 0.5924572112923694, 15)
    11 = @(0.2381956871565345, 0.3847249399981075, 0.1907396565338637, 0.4367375144567543, 0.3798820258321563, 14)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}

$wb.Save()
